# Update countries & provincias Spain
# Applies the diff between before.xlsx and the edited "paises.xlsx":
#  - refreshed COVID-19 counters for several countries (row data on sheet "Pais")
#  - Belice / Nueva Caledonia swap rows (shared-string table reorder upstream)
#  - "Datos actualizados" timestamp bumped from 09:34 to 10:04

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp (A1) -------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 9 de Mayo de 2020 a las 10:04"

# --- Row 8: Rusia -----------------------------------------------------------
$ws.Cells.Item(8, 2).Value = 198676   # Casos totales
$ws.Cells.Item(8, 3).Value = 10817    # Nuevos casos
$ws.Cells.Item(8, 4).Value = 31916    # Casos activos
$ws.Cells.Item(8, 5).Value = 164933   # Recuperados
$ws.Cells.Item(8, 7).Value = 104      # Muertes hoy
$ws.Cells.Item(8, 8).Value = 1827     # Muertes

# --- Row 36: Polonia ---------------------------------------------------------
$ws.Cells.Item(36, 4).Value = 5437    # Casos activos
$ws.Cells.Item(36, 5).Value = 9153    # Recuperados

# --- Row 49: Chequia ---------------------------------------------------------
$ws.Cells.Item(49, 6).Value = 47      # Casos criticos

# --- Row 83: Estonia ----------------------------------------------------------
$ws.Cells.Item(83, 2).Value = 1733    # Casos totales
$ws.Cells.Item(83, 3).Value = 8       # Nuevos casos
$ws.Cells.Item(83, 4).Value = 747     # Casos activos
$ws.Cells.Item(83, 5).Value = 926     # Recuperados
$ws.Cells.Item(83, 6).Value = 5       # Casos criticos
$ws.Cells.Item(83, 7).Value = 4       # Muertes hoy
$ws.Cells.Item(83, 8).Value = 60      # Muertes

# --- Row 88: Eslovaquia -------------------------------------------------------
$ws.Cells.Item(88, 4).Value = 919     # Casos activos
$ws.Cells.Item(88, 5).Value = 510     # Recuperados

# --- Row 103: Sri Lanka -------------------------------------------------------
$ws.Cells.Item(103, 4).Value = 255    # Casos activos
$ws.Cells.Item(103, 5).Value = 571    # Recuperados

# --- Rows 192/193: Nueva Caledonia / Belice swap places -----------------------
# Before: row192 = Nueva Caledonia, row193 = Belice.
# After:  row192 = Belice,          row193 = Nueva Caledonia.
for ($col = 1; $col -le 8; $col++) {
    $v192 = $ws.Cells.Item(192, $col).Value2
    $v193 = $ws.Cells.Item(193, $col).Value2
    $ws.Cells.Item(192, $col).Value = $v193
    $ws.Cells.Item(193, $col).Value = $v192
}
